$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.308.15"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.556.93"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "606.49"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "144.41"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "3.556.69"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").Value = "7.83"
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "4.162.51"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "30.24"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "3.560.94"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "66.411.21"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "11.50"
$ws.Range("E19").Value = "  +5.19%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "14.82"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "431.40"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").Value = "0.610"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").Value = "79.53"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "3.700.23"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "9.15"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "7.95"
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "25.44"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "3.553.25"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E35").Value = "  -6.32%  "
$ws.Range("D36").Value = "7.84"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "5.60"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "176.18"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "5.19"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").Value = "0.888"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").Value = "45.95"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "25.21"
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("D50").Value = "7.14"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "23.19"
$ws.Range("E51").Value = "  +2.49%  "
